# Merge with dev branch: add new CIMSOURCE_* GraphQL test rows to the
# "getDataGraphQL" sheet (rows 12-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getDataGraphQL")

$rows = @(
    @("iems-config-CIMSOURCE_TRAINCONFIG",        '{CIMSOURCE_TRAINCONFIG(cond:"",order:"id DESC") { createTime id optimize_target_type_id simu_id train_cfg_ANN_YN train_cfg_LR_YN train_cfg_RFR_YN train_cfg_SVR_YN train_cfg_SmpEndTime train_cfg_SmpStTime train_cfg_Tree_YN train_cfg_aDa_YN train_cfg_obj_PT train_cfg_obligate1 train_cfg_obligate2 train_cfg_obligate3 train_cfg_obligate4 train_cfg_obligate5 train_cfg_obligate6 train_cfg_predict train_cfg_sample_PT train_cfg_time train_cfg_timetrg_YN}}'),
    @("iems-config-CIMSOURCE_TRAINRESULT",         '{CIMSOURCE_TRAINRESULT(cond:"",order:"id ASC") { variables train_result_obligate4 train_result_obligate5 train_result_obligate2 train_result_obligate3 optimize_target_type_id simu_id train_result_obligate1 config_id train_result_obligate6 id runid }}'),
    @("iems-config-CIMSOURCE_PREDICTCONFIG",       '{CIMSOURCE_PREDICTCONFIG(cond:"",order:"id ASC") { predict_train_model createTime predict_train_obj_wgt_defval predict_cfg_obligate6 predict_cfg_obligate5 predict_cfg_obligate4 predict_cfg_obligate3 predict_cfg_obligate2 predict_cfg_obligate1 optimize_target_type_id simu_id predict_cfg_period predict_cfg_day predict_cfg_time id } }'),
    @("iems-config-CIMSOURCE_PREDICTRESULT",       '{CIMSOURCE_PREDICTRESULT(cond:"",order:"config_id ASC") { Predict_result_time Predict_result_aDa Predict_result_ANN optimize_target_type_id Predict_result_SP5 Predict_result_Tree Predict_result_LR Predict_result_SVR runid id createTime variables predict_train_obj_wgt_defval Predict_result_obligate2 Predict_result_obligate1 Predict_result_RFR Predict_result_obligate6 simu_id Predict_result_obligate5 Predict_result_obligate4 Predict_result_obligate3 config_id Predict_result_AVG } }'),
    @("iems-config-CIMSOURCE_OPTIMALIZERESULT",    '{CIMSOURCE_PREDICTRESULT(cond:"",order:"config_id ASC") { Predict_result_time Predict_result_aDa Predict_result_ANN optimize_target_type_id Predict_result_SP5 Predict_result_Tree Predict_result_LR Predict_result_SVR runid id createTime variables predict_train_obj_wgt_defval Predict_result_obligate2 Predict_result_obligate1 Predict_result_RFR Predict_result_obligate6 simu_id Predict_result_obligate5 Predict_result_obligate4 Predict_result_obligate3 config_id Predict_result_AVG } }'),
    @("iems-config-CIMSOURCE_OPTIMALIZECONFIG",    '{CIMSOURCE_OPTIMALIZECONFIG(cond:"",order:"id ASC") { isActive optimize_target_type_id optimalizeType opt_cfg_opt_starttime id opt_cfg_user_def opt_pred_model opt_cfg_obligate3 rollTasksId opt_cfg_obligate2 opt_cfg_obligate5 createTime opt_cfg_obligate4 opt_cfg_obligate6 opt_cfg_emi_frac opt_name dash_cfg_id simu_id opt_cfg_opt_day opt_cfg_opt_frequency isRTCOP opt_cfg_cost_frac opt_methods rollType opt_cfg_obligate1 } }'),
    @("iems-config-CIMSOURCE_SIMULATION",          '{CIMSOURCE_SIMULATION(cond:"",order:"id ASC") { updatatime historymap chartsetting createtime rtcountmap hiscountmap isTraining obligate1 obligate2 channels defcountmap user_id ispredict name id obligate5 obligate6 treefile obligate3 obligate4 } }'),
    @("iems-config-CIMSOURCE_TOPOLOGY",            '{CIMSOURCE_TOPOLOGY(cond:"",order:"id ASC") { create_time topo_name simu_id obligate1 obligate2 topo_xml update_time update_user parent_id id obligate5 obligate6 obligate3 obligate4 } }'),
    @("iems-config-CIMSOURCE_OPTIMIZE_TARGET_TYPE",'{CIMSOURCE_OPTIMIZE_TARGET_TYPE(cond:"",order:"id ASC") { obligate1 obligate2 optimize_topo_selection optimize_type optimize_name id obligate5 optimize_data_type obligate3 obligate4 } }'),
    @("iems-config-CIMSOURCE_DASHBOARDCONFIG",     '{CIMSOURCE_DASHBOARDCONFIG(cond:"",order:"id ASC") { dashparam name id simu_id optimize_target_type_id } }'),
    @("iems-config-CIMSOURCE_METADATA",            '{CIMSOURCE_METADATA(cond:"",order:"id ASC") { obligate1 obligate2 name id obligate5 params category obligate3 obligate4 } }'),
    @("iems-config-CIMSOURCE_SECTIONALIZEDLINEARIZATION", '{CIMSOURCE_SECTIONALIZEDLINEARIZATION(cond:"",order:"id ASC") { tensor_index_counts target obligate1 obligate2 device_name model_name param id obligate5 obligate3 obligate4 tensor_data } }')
)

$startRow = 12
$endRow = $startRow + $rows.Count - 1
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = "good request, data retrieved"
    $ws.Cells.Item($r, 3).Value = $rows[$i][1]
    $ws.Cells.Item($r, 4).Value = 200
    $ws.Cells.Item($r, 5).Value = 100000
    $ws.Cells.Item($r, 6).Value = "Successfully"
}

# Match the bordered look of the existing data rows for columns B and D
# (columns A, C, E, F already pick up the correct look from the column
# default style).
$ws.Range("B2:B11").Copy()
$ws.Range("B$startRow`:B$endRow").PasteSpecial(-4122)
$ws.Range("D2:D11").Copy()
$ws.Range("D$startRow`:D$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Bump the zoom level on both sheets (matches the refreshed session/view
# state saved with this workbook) and leave the selection where the
# editor left it on the new last row of data.
$ws1 = $wb.Worksheets.Item("getDataEntities")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 130

$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("C25").Select() | Out-Null
